$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values for rows 3-28 (schedule shifted by one week, names reordered)
$ws.Cells.Item(3, 2).Value = "July9th"
$ws.Cells.Item(3, 3).Value = "July10th"
$ws.Cells.Item(3, 4).Value = "July11th"
$ws.Cells.Item(3, 5).Value = "July12th"
$ws.Cells.Item(3, 6).Value = "July13th"
$ws.Cells.Item(3, 7).Value = "July14th"
$ws.Cells.Item(3, 8).Value = "July15th"
$ws.Cells.Item(5, 4).Value = "3:30-8"
$ws.Cells.Item(5, 5).Value = "Meet"
$ws.Cells.Item(5, 7).Value = "3:30-8"
$ws.Cells.Item(6, 2).Value = "3:30-8"
$ws.Cells.Item(6, 4).Value = "3:30-8"
$ws.Cells.Item(6, 5).Value = "1:00-6:00"
$ws.Cells.Item(6, 6).Value = "X"
$ws.Cells.Item(6, 7).Value = "3:30-8"
$ws.Cells.Item(7, 2).Value = "X"
$ws.Cells.Item(7, 3).Value = "X"
$ws.Cells.Item(7, 4).Value = "3:30-8"
$ws.Cells.Item(7, 5).Value = "OFF"
$ws.Cells.Item(7, 6).Value = "10:00-3:30"
$ws.Cells.Item(7, 7).Value = "1:00-6:00"
$ws.Cells.Item(8, 3).Value = "10:15-3:30"
$ws.Cells.Item(8, 4).Value = "3:30-8"
$ws.Cells.Item(8, 5).Value = "Meet"
$ws.Cells.Item(8, 6).Value = "X"
$ws.Cells.Item(8, 7).Value = "X"
$ws.Cells.Item(9, 2).Value = "10:15-3:30"
$ws.Cells.Item(9, 3).Value = "3:30-8"
$ws.Cells.Item(9, 4).Value = "OFF"
$ws.Cells.Item(9, 6).Value = "1:00-6:00"
$ws.Cells.Item(10, 3).Value = "X"
$ws.Cells.Item(10, 4).Value = "X"
$ws.Cells.Item(10, 5).Value = "Meet"
$ws.Cells.Item(11, 3).Value = "10:30-3:30"
$ws.Cells.Item(11, 4).Value = "OFF"
$ws.Cells.Item(11, 5).Value = "4:00-9"
$ws.Cells.Item(11, 6).Value = "10:00-3:30"
$ws.Cells.Item(12, 2).Value = "X"
$ws.Cells.Item(12, 3).Value = "3:30-8"
$ws.Cells.Item(12, 4).Value = "10:15-3:30"
$ws.Cells.Item(12, 5).Value = "4:00-9"
$ws.Cells.Item(13, 2).Value = "3:30-8"
$ws.Cells.Item(13, 3).Value = "OFF"
$ws.Cells.Item(13, 4).Value = "OFF"
$ws.Cells.Item(13, 5).Value = "10:15-4"
$ws.Cells.Item(13, 6).Value = "3:30-8"
$ws.Cells.Item(13, 7).Value = "10:30-3:30"
$ws.Cells.Item(14, 4).Value = "OFF"
$ws.Cells.Item(15, 2).Value = "X"
$ws.Cells.Item(15, 3).Value = "OFF"
$ws.Cells.Item(15, 4).Value = "10:15-3:30"
$ws.Cells.Item(15, 7).Value = "10:45-3:30"
$ws.Cells.Item(16, 1).Value = "Nathan Debergh"
$ws.Cells.Item(16, 2).Value = "X"
$ws.Cells.Item(16, 3).Value = "X"
$ws.Cells.Item(16, 4).Value = "X"
$ws.Cells.Item(16, 5).Value = "X"
$ws.Cells.Item(16, 6).Value = "X"
$ws.Cells.Item(16, 7).Value = "X"
$ws.Cells.Item(17, 1).Value = "Phillip Thompson"
$ws.Cells.Item(17, 2).Value = "OFF"
$ws.Cells.Item(17, 3).Value = "3:30-8"
$ws.Cells.Item(17, 4).Value = "OFF"
$ws.Cells.Item(17, 5).Value = "4:00-9"
$ws.Cells.Item(17, 6).Value = "10:15-3:30"
$ws.Cells.Item(17, 7).Value = "OFF"
$ws.Cells.Item(18, 1).Value = "Madison Johnson"
$ws.Cells.Item(18, 2).Value = "X"
$ws.Cells.Item(19, 1).Value = "Asher Bobbett"
$ws.Cells.Item(19, 2).Value = "10:15-3:30"
$ws.Cells.Item(19, 3).Value = "10:30-3:30"
$ws.Cells.Item(19, 4).Value = "10:15-3:30"
$ws.Cells.Item(19, 5).Value = "OFF"
$ws.Cells.Item(19, 7).Value = "X"
$ws.Cells.Item(20, 1).Value = "Blake Ucherek"
$ws.Cells.Item(20, 3).Value = "OFF"
$ws.Cells.Item(20, 4).Value = "3:30-8"
$ws.Cells.Item(20, 5).Value = "10:15-4"
$ws.Cells.Item(20, 6).Value = "3:30-8"
$ws.Cells.Item(20, 7).Value = "OFF"
$ws.Cells.Item(21, 1).Value = "Ethan Van Horn "
$ws.Cells.Item(21, 2).Value = "OFF"
$ws.Cells.Item(21, 4).Value = "OFF"
$ws.Cells.Item(21, 5).Value = "10:30-4"
$ws.Cells.Item(21, 6).Value = "10:15-3:30"
$ws.Cells.Item(22, 1).Value = "Kai King"
$ws.Cells.Item(22, 2).Value = "OFF"
$ws.Cells.Item(22, 3).Value = "OFF"
$ws.Cells.Item(22, 4).Value = "10:30-3:30"
$ws.Cells.Item(22, 5).Value = "Meet"
$ws.Cells.Item(22, 7).Value = "3:30-8"
$ws.Cells.Item(23, 1).Value = "Madeline Ellison"
$ws.Cells.Item(23, 2).Value = "OFF"
$ws.Cells.Item(23, 3).Value = "OFF"
$ws.Cells.Item(23, 4).Value = "3:30-8"
$ws.Cells.Item(23, 5).Value = "OFF"
$ws.Cells.Item(23, 6).Value = "3:30-8"
$ws.Cells.Item(23, 7).Value = "10:45-3:30"
$ws.Cells.Item(24, 1).Value = "Tyler Carpenter"
$ws.Cells.Item(24, 2).Value = "10:30-3:30"
$ws.Cells.Item(24, 3).Value = "3:30-8"
$ws.Cells.Item(24, 4).Value = "OFF"
$ws.Cells.Item(24, 5).Value = "10:30-4"
$ws.Cells.Item(25, 2).Value = "10:30-3:30"
$ws.Cells.Item(25, 4).Value = "10:30-3:30"
$ws.Cells.Item(25, 5).Value = "Meet"
$ws.Cells.Item(25, 6).Value = "10:15-3:30"
$ws.Cells.Item(25, 7).Value = "OFF"
$ws.Cells.Item(26, 2).Value = "OFF"
$ws.Cells.Item(26, 3).Value = "3:30-8"
$ws.Cells.Item(26, 4).Value = "10:30-3:30"
$ws.Cells.Item(26, 6).Value = "3:30-8"
$ws.Cells.Item(26, 7).Value = "3:30-8"
$ws.Cells.Item(27, 2).Value = "3:30-8"
$ws.Cells.Item(27, 3).Value = "OFF"
$ws.Cells.Item(27, 4).Value = "OFF"
$ws.Cells.Item(27, 5).Value = "OFF"
$ws.Cells.Item(28, 1).Value = "Brent Horwitz"
$ws.Cells.Item(28, 2).Value = "3:30-8"
$ws.Cells.Item(28, 3).Value = "10:30-3:30"
$ws.Cells.Item(28, 4).Value = "OFF"
$ws.Cells.Item(28, 5).Value = "4:00-9"
$ws.Cells.Item(28, 6).Value = "X"
$ws.Cells.Item(28, 7).Value = "X"

# Remove the last row (Brent Horwitz row 29) - data shifted up, sheet now ends at row 28
$ws.Rows.Item(29).Delete()

Write-Host "Schedule updated successfully"